$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4916.6665
$ws.Range("I43").Value = 5260
$ws.Range("K43").Value = 5260
$ws.Range("M43").Value = -5191
$ws.Range("H64").Value = 12220.053
$ws.Range("I64").Value = 6329.25
$ws.Range("J64").Value = 13790.934
$ws.Range("K64").Value = 6329.25
$ws.Range("L64").Value = 13790.934
$ws.Range("M64").Value = -6081.25
$ws.Range("N64").Value = -14286.934
$ws.Range("H67").Value = 12220.053
$ws.Range("I67").Value = 6329.25
$ws.Range("J67").Value = 13790.934
$ws.Range("K67").Value = 6329.25
$ws.Range("L67").Value = 13790.934
$ws.Range("M67").Value = -5471.25
$ws.Range("N67").Value = -15506.934
$ws.Range("H94").Value = 925.3570999999999
$ws.Range("I94").Value = 925.3570999999999
$ws.Range("K94").Value = 925.3570999999999
$ws.Range("M94").Value = -474.3570999999999
$ws.Range("H116").Value = 8553.462
$ws.Range("I116").Value = 5840
$ws.Range("J116").Value = 10249.375
$ws.Range("K116").Value = 5840
$ws.Range("L116").Value = 10249.375
$ws.Range("M116").Value = -2398
$ws.Range("N116").Value = -17133.375
$ws.Range("H137").Value = 5402.6
$ws.Range("I137").Value = 6690.6113
$ws.Range("K137").Value = 20071.8339
$ws.Range("M137").Value = -17521.8339
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 120110.2
$ws.Range("I34").Value = 10000
$ws.Range("J34").Value = 147637.75
$ws.Range("K34").Value = 10000
$ws.Range("L34").Value = 147637.75
$ws.Range("M34").Value = -9729
$ws.Range("N34").Value = -148179.75
$ws.Range("H45").Value = 2456.739
$ws.Range("I45").Value = 2026.579
$ws.Range("J45").Value = 4500
$ws.Range("K45").Value = 2026.579
$ws.Range("L45").Value = 4500
$ws.Range("M45").Value = -1649.579
$ws.Range("N45").Value = -5254
$ws.Range("H102").Value = 2622.7778
$ws.Range("I102").Value = 2825.625
$ws.Range("K102").Value = 2825.625
$ws.Range("M102").Value = -1203.625
$ws.Range("H132").Value = 1834.6316
$ws.Range("I132").Value = 1659.8889
$ws.Range("K132").Value = 4979.6667
$ws.Range("M132").Value = -2449.6667
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 544.5294
$ws.Range("I80").Value = 399.75
$ws.Range("J80").Value = 589.0769
$ws.Range("K80").Value = 399.75
$ws.Range("L80").Value = 589.0769
$ws.Range("M80").Value = 598.25
$ws.Range("N80").Value = -2585.0769
$ws.Range("H82").Value = 4199
$ws.Range("I82").Value = 4199
$ws.Range("K82").Value = 4199
$ws.Range("M82").Value = -3816
$ws.Range("H83").Value = 544.5294
$ws.Range("I83").Value = 399.75
$ws.Range("J83").Value = 589.0769
$ws.Range("K83").Value = 1998.75
$ws.Range("L83").Value = 2945.3845
$ws.Range("M83").Value = 2993.25
$ws.Range("N83").Value = -12929.3845
$ws.Range("H85").Value = 4199
$ws.Range("I85").Value = 4199
$ws.Range("K85").Value = 4199
$ws.Range("M85").Value = -2873
$ws.Range("H134").Value = 1968.6666
$ws.Range("I134").Value = 1968.6666
$ws.Range("K134").Value = 5905.9998
$ws.Range("M134").Value = -3370.9998
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1282.2222
$ws.Range("I99").Value = 1282.2222
$ws.Range("K99").Value = 1282.2222
$ws.Range("M99").Value = 215.7778000000001
$ws.Range("H107").Value = 5347.909
$ws.Range("I107").Value = 584
$ws.Range("K107").Value = 584
$ws.Range("M107").Value = 1336
$ws.Range("H126").Value = 1282.2222
$ws.Range("I126").Value = 1282.2222
$ws.Range("K126").Value = 3846.6666
$ws.Range("M126").Value = -1376.6666
$ws.Range("H134").Value = 2345.125
$ws.Range("I134").Value = 1696.1428
$ws.Range("K134").Value = 5088.428400000001
$ws.Range("M134").Value = -2553.428400000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2281.9722
$ws.Range("I107").Value = 2986.1428
$ws.Range("J107").Value = 2112
$ws.Range("K107").Value = 8958.428400000001
$ws.Range("L107").Value = 6336
$ws.Range("M107").Value = -7038.428400000001
$ws.Range("N107").Value = -10176
$ws.Range("H113").Value = 1254.25
$ws.Range("I113").Value = 1435
$ws.Range("J113").Value = 1073.5
$ws.Range("K113").Value = 4305
$ws.Range("L113").Value = 3220.5
$ws.Range("M113").Value = -2135
$ws.Range("N113").Value = -7560.5
$ws.Range("H140").Value = 2772.7273
$ws.Range("I140").Value = 2638.75
$ws.Range("K140").Value = 7916.25
$ws.Range("M140").Value = -2736.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 59735.65
$ws.Range("I80").Value = 141651.38
$ws.Range("J80").Value = 5125.1665
$ws.Range("K80").Value = 141651.38
$ws.Range("L80").Value = 5125.1665
$ws.Range("M80").Value = -140653.38
$ws.Range("N80").Value = -7121.1665
$ws.Range("H83").Value = 59735.65
$ws.Range("I83").Value = 141651.38
$ws.Range("J83").Value = 5125.1665
$ws.Range("K83").Value = 708256.9
$ws.Range("L83").Value = 25625.8325
$ws.Range("M83").Value = -703264.9
$ws.Range("N83").Value = -35609.8325
$ws.Range("H97").Value = 5714.684
$ws.Range("J97").Value = 16969.834
$ws.Range("L97").Value = 16969.834
$ws.Range("N97").Value = -17961.834
$ws.Range("H113").Value = 5085.067
$ws.Range("I113").Value = 3261.5
$ws.Range("J113").Value = 7169.143
$ws.Range("K113").Value = 3261.5
$ws.Range("L113").Value = 7169.143
$ws.Range("M113").Value = -1091.5
$ws.Range("N113").Value = -11509.143
$ws.Range("H122").Value = 3906.9333
$ws.Range("I122").Value = 3466
$ws.Range("J122").Value = 4568.3335
$ws.Range("K122").Value = 10398
$ws.Range("L122").Value = 13705.0005
$ws.Range("M122").Value = -7948
$ws.Range("N122").Value = -18605.0005
$ws.Range("H135").Value = 39647.06
$ws.Range("J135").Value = 39647.06
$ws.Range("L135").Value = 39647.06
$ws.Range("N135").Value = -49787.06
$ws.Range("H140").Value = 110000
$ws.Range("J140").Value = 110000
$ws.Range("L140").Value = 110000
$ws.Range("N140").Value = -120360
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10101989
$ws.Range("I22").Value = 22727694
$ws.Range("J22").Value = 1424.6
$ws.Range("K22").Value = 22727694
$ws.Range("L22").Value = 1424.6
$ws.Range("M22").Value = -22727399
$ws.Range("N22").Value = -2014.6
$ws.Range("H27").Value = 10101989
$ws.Range("I27").Value = 22727694
$ws.Range("J27").Value = 1424.6
$ws.Range("K27").Value = 22727694
$ws.Range("L27").Value = 1424.6
$ws.Range("M27").Value = -22727587
$ws.Range("N27").Value = -1638.6
$ws.Range("H40").Value = 30063.125
$ws.Range("I40").Value = 4250.8335
$ws.Range("K40").Value = 4250.8335
$ws.Range("M40").Value = -4114.8335
$ws.Range("H93").Value = 2072.1428
$ws.Range("I93").Value = 1917.6666
$ws.Range("K93").Value = 1917.6666
$ws.Range("M93").Value = -669.6666
$ws.Range("H136").Value = 2694.276
$ws.Range("I136").Value = 2074.4211
$ws.Range("J136").Value = 3872
$ws.Range("K136").Value = 6223.263300000001
$ws.Range("L136").Value = 11616
$ws.Range("M136").Value = -3673.263300000001
$ws.Range("N136").Value = -16716
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 424.1
$ws.Range("I107").Value = 356.5
$ws.Range("J107").Value = 525.5
$ws.Range("K107").Value = 1069.5
$ws.Range("L107").Value = 1576.5
$ws.Range("M107").Value = 850.5
$ws.Range("N107").Value = -5416.5
$ws.Range("H115").Value = 29833.334
$ws.Range("J115").Value = 29833.334
$ws.Range("L115").Value = 29833.334
$ws.Range("N115").Value = -32967.334
$ws.Range("H132").Value = 3749.4546
$ws.Range("I132").Value = 3244.8667
$ws.Range("K132").Value = 9734.6001
$ws.Range("M132").Value = -7204.6001
